$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.231869697570801
$ws.Range("B1").Value = 2.400803804397583
$ws.Range("C1").Value = 1.995851159095764
$ws.Range("D1").Value = 1.896677613258362
$ws.Range("E1").Value = 1.684993863105774
